$d = $word.ActiveDocument

# The run currently reads:
#   " observation is grouped into a cluster it is stuck there). The final algorithm
#    calculates the p-values for the hierarchical clustering based on multiscale
#    bootstrap resampling. This will help determine if the clusters calculated are
#    supported by the data. This helps with understanding the graph produced by
#    RStudio."
#
# Split it into three paragraphs:
#   1) "... it is stuck there). "                                   (unchanged paragraph/style)
#   2) "Research Model-Based Clustering"                             (new Heading 3)
#   3) "The Research Model-Based clustering algorithm calculates the p-values ..." (new body)
$old = "The final algorithm calculates the p-values for the hierarchical clustering based on multiscale bootstrap resampling. This will help determine if the clusters calculated are supported by the data. This helps with understanding the graph produced by RStudio."
$new = "^pResearch Model-Based Clustering^pThe Research Model-Based clustering algorithm calculates the p-values for the hierarchical clustering based on multiscale bootstrap resampling. This will help determine if the clusters calculated are supported by the data. This helps with understanding the graph produced by RStudio."

$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# Locate the freshly split heading/body paragraphs by scanning for their text
# (robust to any paragraph-index drift caused by the Find/Replace).
$headingPara = $null
$bodyPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Research Model-Based Clustering`r") {
        $headingPara = $p
        $bodyPara = $d.Paragraphs($i + 1)
        break
    }
}

# New heading paragraph -> Heading 3 style (matches the other section headings, e.g.
# "Hierarchical Clustering" just above it).
$headingPara.Style = "Heading 3"

# New body paragraph -> split into the three runs used by the target markup:
#   "The " | "Research Model-Based clustering" | " algorithm calculates ... RStudio."
# and make sure the whole paragraph carries the Times New Roman font (ascii/hAnsi/cs)
# used throughout the rest of the document.
$bStart = $bodyPara.Range.Start
$bEnd = $bodyPara.Range.End - 1   # exclude the paragraph mark

$run1 = $d.Range($bStart, $bStart + 4)                 # "The "
$run2 = $d.Range($bStart + 4, $bStart + 4 + 31)         # "Research Model-Based clustering"
$run3 = $d.Range($bStart + 4 + 31, $bEnd)               # " algorithm calculates ... RStudio."

$run1.Font.Name = "Times New Roman"
$run2.Font.Name = "Times New Roman"
$run3.Font.Name = "Times New Roman"

# Fill in the complex-script (w:cs) font across the whole paragraph in one shot --
# setting NameBi on the freshly-split sub-ranges individually doesn't stick.
$bodyPara.Range.Font.NameBi = "Times New Roman"
